$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 285, pushing existing rows 285..395 down to 286..396.
$ws.Rows("285:285").Insert()

# Populate the newly inserted row 285 with the new weekly price record.
$ws.Cells.Item(285, 1).Value = 10
$ws.Cells.Item(285, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(285, 3).Value = "La Araucanía"
$ws.Cells.Item(285, 4).Value = 44924
$ws.Cells.Item(285, 5).Value = 9
$ws.Cells.Item(285, 6).Value = 100112017
$ws.Cells.Item(285, 7).Value = "Apio"
$ws.Cells.Item(285, 8).Value = "Americana (o)"
$ws.Cells.Item(285, 9).Value = "Primera"
$ws.Cells.Item(285, 10).Value = 85
$ws.Cells.Item(285, 11).Value = 10000
$ws.Cells.Item(285, 12).Value = 10000
$ws.Cells.Item(285, 13).Value = 10000
$ws.Cells.Item(285, 14).Value = "$/docena de matas"
$ws.Cells.Item(285, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(285, 16).Value = 1667
$ws.Cells.Item(285, 17).Value = 6
$ws.Cells.Item(285, 18).Value = "Hortaliza"
